$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: Problem #3 -> "Is it submitted before the deadline?" answer
# flips from "Yes" to "No" (bold), as its own new run (the diff removes
# <w:t>Yes</w:t> from the tab-bearing run and appends a fresh bold run
# carrying "No").
# ------------------------------------------------------------------
$p3Start = -1
$p4Start = -1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Problem #3*") { $p3Start = $p.Range.Start }
    if ($p.Range.Text -like "Problem #4*") { $p4Start = $p.Range.Start }
}

$deadlinePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -gt $p3Start -and $p.Range.Start -lt $p4Start -and `
        $p.Range.Text -like "Is it submitted before the deadline?*") {
        $deadlinePara = $p
    }
}

$yesRng = $deadlinePara.Range.Duplicate
$yesRng.Find.Execute("Yes", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($yesRng.Find.Found) {
    $insertAt = $yesRng.Start
    $yesRng.Delete()
    $noRng = $d.Range($insertAt, $insertAt)
    $noRng.InsertAfter("No")
    $noFmt = $d.Range($insertAt, $insertAt + 2)
    $noFmt.Bold = 1
}

# ------------------------------------------------------------------
# Change 2: Problem #4 -> "Does your code run without any compilation
# errors?" answer was split across three runs; collapse it back into a
# single bold run with the full sentence.
# ------------------------------------------------------------------
$needle = $d.Content.Duplicate
$needle.Find.Execute("No (something wrong with castings I couldn")
if ($needle.Find.Found) {
    $containingPara = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Start -le $needle.Start -and $p.Range.End -ge $needle.End) {
            $containingPara = $p
        }
    }
    $fullRng = $d.Range($needle.Start, $containingPara.Range.End)
    $newText = "No (something wrong with castings I couldn" + [char]0x2019 + `
        "t figure it out, but I think the implementations are correct)"
    $fullRng.Text = $newText
    $boldRng = $d.Range($needle.Start, $needle.Start + $newText.Length)
    $boldRng.Bold = 1
}
